$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.844.21'
$ws.Range('E2').Value = '  -1.94%  '

# Row 3
$ws.Range('D3').Value = '3.560.43'
$ws.Range('E3').Value = '  -3.57%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '568.68'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -5.49%  '

# Row 6
$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '190.07'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.43%  '

# Row 7
$ws.Range('D7').Value = '3.555.96'
$ws.Range('E7').Value = '  -3.63%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.612'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -2.18%  '

# Row 9
$ws.Range('E9').Value = '  +0.09%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.673'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -5.41%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '55.44'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -4.41%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.148'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -4.74%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000267'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -4.14%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.79'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -4.64%  '

# Row 15
$ws.Range('D15').Value = '4.136.81'
$ws.Range('E15').Value = '  -3.63%  '

# Row 16
$ws.Range('D16').Value = '3.569.87'
$ws.Range('E16').Value = '  -3.59%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.125'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.39%  '

# Row 18
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '66.806.18'
$ws.Range('E18').Value = '  -1.95%  '

# Row 19
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '18.20'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -5.00%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.10'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -4.36%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.06'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -6.50%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '398.24'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.24%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.16'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -7.63%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '85.63'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -3.81%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '11.56'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.30%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.91'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.13%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '12.36'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.81%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.09'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.98%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '3.61'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -2.67%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.69'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.86%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.88'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -5.79%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '31.05'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.36%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '637.80'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.50%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '12.05'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.39%  '

# Row 35
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '63.81'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -5.12%  '

# Row 36
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.113'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.77%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '41.94'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -8.83%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.400'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.83%  '

# Row 39
$ws.Range('E39').Value = '  +0.14%  '

# Row 40
$ws.Range('D40').Value = '0.0₃0756'
$ws.Range('E40').Value = '  -5.35%  '

# Row 41
$ws.Range('D41').Value = '3.187.79'
$ws.Range('E41').Value = '  +11.20%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.132'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.81%  '

# Row 43
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.20%  '

# Row 44
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.68'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +4.21%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.93'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.29%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0411'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.81%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.09'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.72%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.129'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -6.07%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '140.98'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.11%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.50'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -6.34%  '

# Row 51
$ws.Range('B51').Value = 'WEMIXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.52'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -5.23%  '
